# Commit: "Se agrega el scrap de Plumbersstock Auto_scrap_Plumbersstock.py."
#
# Concrete edits observed in the OOXML diff:
#   1. The (only) worksheet is renamed "HomeDepot_URL" -> "URL".
#   2. Cell G4 ("Elongated  " with trailing spaces) is corrected to the
#      already-existing shared string "Elongated" (no trailing spaces).
#      Dropping the now-unused "Elongated  " shared string is what makes
#      the sharedStrings table shrink and every later shared-string index
#      shift down by one in the diff - that is a side effect of this one
#      value fix, not a separate edit.
#   3. The sheet's remembered selection moves to C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "URL"

$ws.Range("G4").Value = "Elongated"

$ws.Range("C31").Select()
